$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 22:52"

# Row 21 -> Asturias (new data pushed to top of this block)
$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 2170
$ws.Range("C21").Value = 522
$ws.Range("D21").Value = 1480
$ws.Range("E21").Value = 168

# Row 22 -> Gipuzkoa/Guipuzcoa (shifted down from row 21)
$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B22").Value = 2161
$ws.Range("C22").Value = 5428
$ws.Range("D22").Value = 5145
$ws.Range("E22").Value = 161

# Row 23 -> Sevilla (shifted down from row 22)
$ws.Range("A23").Value = "Sevilla"
$ws.Range("B23").Value = 2159
$ws.Range("C23").Value = 327
$ws.Range("D23").Value = 1650
$ws.Range("E23").Value = 182

# Row 24 -> Segovia (shifted down from row 23)
$ws.Range("A24").Value = "Segovia"
$ws.Range("B24").Value = 2103
$ws.Range("C24").Value = 582
$ws.Range("D24").Value = 1363
$ws.Range("E24").Value = 158

# Row 30 -> Murcia (new data pushed above Pontevedra)
$ws.Range("A30").Value = "Murcia"
$ws.Range("B30").Value = 1598
$ws.Range("C30").Value = 591
$ws.Range("D30").Value = 896
$ws.Range("E30").Value = 111

# Row 31 -> Pontevedra (shifted down from row 30)
$ws.Range("A31").Value = "Pontevedra"
$ws.Range("B31").Value = 1536
$ws.Range("C31").Value = 333
$ws.Range("D31").Value = 1411
$ws.Range("E31").Value = 30

# Row 33 -> Tenerife values updated
$ws.Range("B33").Value = 1248
$ws.Range("C33").Value = 378
$ws.Range("D33").Value = 788
$ws.Range("E33").Value = 72

# Row 50 -> Gran Canaria values updated
$ws.Range("B50").Value = 470
$ws.Range("C50").Value = 203
$ws.Range("D50").Value = 239
$ws.Range("E50").Value = 28

# Row 54 -> Melilla values updated
$ws.Range("B54").Value = 103
$ws.Range("C54").Value = 27
$ws.Range("D54").Value = 74

# Row 55 -> Ceuta values updated
$ws.Range("B55").Value = 99
$ws.Range("C55").Value = 39
$ws.Range("D55").Value = 56

# Row 56 -> La Palma, Muertes updated
$ws.Range("E56").Value = 4

# Row 57 -> Lanzarote values updated
$ws.Range("B57").Value = 79
$ws.Range("D57").Value = 57
